$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has one data row (row 2) describing the
# "preproduccion" environment. We need to:
#   1. Push that existing row down to row 3, keeping its formatting/styles
#      and its hyperlink (pointing the hyperlink at the new row).
#   2. Update the (now relocated) NroSiniestro value in row 3.
#   3. Insert a brand-new environment ("ssurgwsoadev4-oci...") as the new
#      row 2.
# ---------------------------------------------------------------------------

# 1) Duplicate row 2's values+formats down into row 3 (keeps shared-string
#    text formatting, including the quote-prefixed NroSiniestro cell).
$ws.Range("A2:D2").Copy()
$ws.Range("A3").PasteSpecial()

$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial()
$ws.Range("B3").Style = $ws.Range("B2").Style

# 2) Move the hyperlink itself from B2 to B3 (Add() re-applies the
#    hyperlink font the first time, so restore the original style right
#    after).
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B3").Style = $ws.Range("B2").Style

# 3) The claim number tied to the preproduccion environment changes now
#    that it lives on row 3.
$ws.Range("E3").Value2 = "1120194100378"

# 4) Overwrite row 2 with the brand-new environment's data.
$ws.Range("A2").Value2 = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value2 = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("C2").Value2 = "dgariffo"
$ws.Range("D2").Value2 = "silverarrow"
$ws.Range("E2").Value2 = "1120194100385"

# 5) Selection / active cell, as recorded by the author while editing.
$ws.Range("L10").Select()
